# Applies the documented edit: completes the Slot Machines algorithm outline
# (INPUT / PROCESS / OUTPUT sections) and relocates the stray "_GoBack" bookmark
# from the title paragraph into the PROCESS section's "Add 9 quarters" bullet.

$d = $word.ActiveDocument

function Replace-FirstLineWith([string]$bodyXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute("first line", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find 'first line' placeholder to replace."
    }
    $target = $d.Range($rng.Start, $rng.End)
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
        $bodyXml + `
        '<w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($pkg)
}

# 1) Drop the leftover "_GoBack" bookmark from the title paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) INPUT section.
$inputXml = '<w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/></w:rPr><w:t xml:space="preserve">Prompt for number of </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/></w:rPr><w:t>quarters</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/></w:rPr><w:t>Ensure (0, 1000)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/></w:rPr><w:t>Store value</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/></w:rPr><w:t xml:space="preserve">Prompt for </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/></w:rPr><w:t># of times first machine was played</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/></w:rPr><w:t>Ensure [</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/></w:rPr><w:t>0, 33]</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/></w:rPr><w:t>Store value</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/></w:rPr><w:t>Prompt for # of times second machine was played</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/></w:rPr><w:t>Ensure [</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/></w:rPr><w:t>0, 98</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/></w:rPr><w:t>]</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Store value</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Prompt for # of times third machine played</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Ensure [</w:t></w:r><w:r><w:t>0, 8]</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Store value</w:t></w:r></w:p>'
Replace-FirstLineWith $inputXml

# 3) PROCESS section.
$processXml = '<w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>loop while quarters &gt; 100</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Check if # of quarters &gt; 0</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>break</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Add one “play” to first machine, subtract one quarter</w:t></w:r><w:r><w:t>, add one “play” to total plays</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>If (</w:t></w:r><w:r><w:t xml:space="preserve">machine #1 </w:t></w:r><w:r><w:t>plays since last win</w:t></w:r><w:r><w:t xml:space="preserve"> ==35</w:t></w:r><w:r><w:t xml:space="preserve">) </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Add </w:t></w:r><w:r><w:t>30 quarters</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Set machine #1 plays since last win to 0</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Check if # of quarters &gt; 0</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Break</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Add one “play” to second machine, subtract one quarter, add one “play” to total plays</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">If (machine #2 plays since last win == </w:t></w:r><w:r><w:t>100)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Add 60 quarters</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Set machine #2 plays since last win to 0</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Check if # of quarters &gt; 0</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>break</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Add one “play” to third machine”, subtract 1 quarter, add one “play” to total plays</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>If (</w:t></w:r><w:r><w:t>machine #3</w:t></w:r><w:r><w:t xml:space="preserve"> plays since last win == </w:t></w:r><w:r><w:t>10</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Add </w:t></w:r><w:r><w:t>9</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> quarters</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Set </w:t></w:r><w:r><w:t>machine #3</w:t></w:r><w:r><w:t xml:space="preserve"> plays since last win to 0</w:t></w:r></w:p>'
Replace-FirstLineWith $processXml

# 4) OUTPUT section.
$outputXml = '<w:p><w:pPr><w:pStyle w:val="Body"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/></w:rPr><w:t>Print(“Martha plays \(totalPlays)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/></w:rPr><w:t xml:space="preserve"> times before going broke.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/></w:rPr><w:t>”)</w:t></w:r></w:p>'
Replace-FirstLineWith $outputXml
